# Auto-generated Excel COM-interop edit script
# Applies value updates to columns H-N across multiple sheets
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 69.5  # H11: 77 -> 69.5
$ws.Cells.Item(11, 9).Value = 69.5  # I11: 77 -> 69.5
$ws.Cells.Item(11, 11).Value = 69.5  # K11: 77 -> 69.5
$ws.Cells.Item(11, 13).Value = 70.5  # M11: 63 -> 70.5
$ws.Cells.Item(29, 8).Value = 1500  # H29: 0 -> 1500
$ws.Cells.Item(29, 9).Value = 1500  # I29: 0 -> 1500
$ws.Cells.Item(29, 11).Value = 4500  # K29: 0 -> 4500
$ws.Cells.Item(29, 13).Value = -4219  # M29: None -> -4219
$ws.Cells.Item(131, 8).Value = 4666.6665  # H131: 2000 -> 4666.6665
$ws.Cells.Item(131, 9).Value = 4666.6665  # I131: 2000 -> 4666.6665
$ws.Cells.Item(131, 11).Value = 13999.9995  # K131: 6000 -> 13999.9995
$ws.Cells.Item(131, 13).Value = -8959.999500000002  # M131: -960 -> -8959.999500000002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 96  # H5: 149.5 -> 96
$ws.Cells.Item(5, 9).Value = 95.5  # I5: 149.5 -> 95.5
$ws.Cells.Item(5, 10).Value = 100  # J5: 0 -> 100
$ws.Cells.Item(5, 11).Value = 95.5  # K5: 149.5 -> 95.5
$ws.Cells.Item(5, 12).Value = 100  # L5: 0 -> 100
$ws.Cells.Item(5, 13).Value = 16.5  # M5: -37.5 -> 16.5
$ws.Cells.Item(5, 14).Value = -324  # N5: None -> -324

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 96  # H4: 149.5 -> 96
$ws.Cells.Item(4, 9).Value = 95.5  # I4: 149.5 -> 95.5
$ws.Cells.Item(4, 10).Value = 100  # J4: 0 -> 100
$ws.Cells.Item(4, 11).Value = 95.5  # K4: 149.5 -> 95.5
$ws.Cells.Item(4, 12).Value = 100  # L4: 0 -> 100
$ws.Cells.Item(4, 13).Value = 19.5  # M4: -34.5 -> 19.5
$ws.Cells.Item(4, 14).Value = -330  # N4: None -> -330
$ws.Cells.Item(107, 8).Value = 2000  # H107: 0 -> 2000
$ws.Cells.Item(107, 9).Value = 2000  # I107: 0 -> 2000
$ws.Cells.Item(107, 11).Value = 2000  # K107: 0 -> 2000
$ws.Cells.Item(107, 13).Value = -80  # M107: None -> -80

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 156.28572  # H7: 147.85715 -> 156.28572
$ws.Cells.Item(7, 9).Value = 105  # I7: 86.25 -> 105
$ws.Cells.Item(7, 10).Value = 194.75  # J7: 230 -> 194.75
$ws.Cells.Item(7, 11).Value = 105  # K7: 86.25 -> 105
$ws.Cells.Item(7, 12).Value = 194.75  # L7: 230 -> 194.75
$ws.Cells.Item(7, 13).Value = 8  # M7: 26.75 -> 8
$ws.Cells.Item(7, 14).Value = -420.75  # N7: -456 -> -420.75
$ws.Cells.Item(70, 8).Value = 25000  # H70: 0 -> 25000
$ws.Cells.Item(70, 10).Value = 25000  # J70: 0 -> 25000
$ws.Cells.Item(70, 12).Value = 25000  # L70: 0 -> 25000
$ws.Cells.Item(70, 14).Value = -25630  # N70: None -> -25630
$ws.Cells.Item(73, 8).Value = 25000  # H73: 0 -> 25000
$ws.Cells.Item(73, 10).Value = 25000  # J73: 0 -> 25000
$ws.Cells.Item(73, 12).Value = 25000  # L73: 0 -> 25000
$ws.Cells.Item(73, 14).Value = -27184  # N73: None -> -27184
$ws.Cells.Item(82, 8).Value = 50000  # H82: 0 -> 50000
$ws.Cells.Item(82, 10).Value = 50000  # J82: 0 -> 50000
$ws.Cells.Item(82, 12).Value = 50000  # L82: 0 -> 50000
$ws.Cells.Item(82, 14).Value = -50722  # N82: None -> -50722
$ws.Cells.Item(85, 8).Value = 50000  # H85: 0 -> 50000
$ws.Cells.Item(85, 10).Value = 50000  # J85: 0 -> 50000
$ws.Cells.Item(85, 12).Value = 50000  # L85: 0 -> 50000
$ws.Cells.Item(85, 14).Value = -52496  # N85: None -> -52496
$ws.Cells.Item(141, 8).Value = 199998.5  # H141: 103791 -> 199998.5
$ws.Cells.Item(141, 9).Value = 0  # I141: 78499 -> 0
$ws.Cells.Item(141, 10).Value = 199998.5  # J141: 112221.664 -> 199998.5
$ws.Cells.Item(141, 11).Value = 0  # K141: 78499 -> 0
$ws.Cells.Item(141, 12).Value = 199998.5  # L141: 112221.664 -> 199998.5
$ws.Cells.Item(141, 13).ClearContents()  # M141: -73319 -> (removed)
$ws.Cells.Item(141, 14).Value = -210358.5  # N141: -122581.664 -> -210358.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 57.153847  # H6: 58.833332 -> 57.153847
$ws.Cells.Item(6, 10).Value = 38.75  # J6: 39.333332 -> 38.75
$ws.Cells.Item(6, 12).Value = 116.25  # L6: 117.999996 -> 116.25
$ws.Cells.Item(6, 14).Value = -342.25  # N6: -343.999996 -> -342.25
$ws.Cells.Item(7, 8).Value = 366  # H7: 499 -> 366
$ws.Cells.Item(7, 9).Value = 51.5  # I7: 3 -> 51.5
$ws.Cells.Item(7, 11).Value = 154.5  # K7: 9 -> 154.5
$ws.Cells.Item(7, 13).Value = -42.5  # M7: 103 -> -42.5
$ws.Cells.Item(9, 8).Value = 1500  # H9: 0 -> 1500
$ws.Cells.Item(9, 10).Value = 1500  # J9: 0 -> 1500
$ws.Cells.Item(9, 12).Value = 4500  # L9: 0 -> 4500
$ws.Cells.Item(9, 14).Value = -4948  # N9: None -> -4948
$ws.Cells.Item(10, 8).Value = 243.6  # H10: 73 -> 243.6
$ws.Cells.Item(10, 9).Value = 54.75  # I10: 73 -> 54.75
$ws.Cells.Item(10, 10).Value = 999  # J10: 0 -> 999
$ws.Cells.Item(10, 11).Value = 164.25  # K10: 219 -> 164.25
$ws.Cells.Item(10, 12).Value = 2997  # L10: 0 -> 2997
$ws.Cells.Item(10, 13).Value = -25.25  # M10: -80 -> -25.25
$ws.Cells.Item(10, 14).Value = -3275  # N10: None -> -3275
$ws.Cells.Item(16, 8).Value = 200  # H16: 0 -> 200
$ws.Cells.Item(16, 9).Value = 200  # I16: 0 -> 200
$ws.Cells.Item(16, 11).Value = 600  # K16: 0 -> 600
$ws.Cells.Item(16, 13).Value = -427  # M16: None -> -427
$ws.Cells.Item(17, 8).Value = 212.5  # H17: 142 -> 212.5
$ws.Cells.Item(17, 9).Value = 0  # I17: 1 -> 0
$ws.Cells.Item(17, 11).Value = 0  # K17: 3 -> 0
$ws.Cells.Item(17, 13).ClearContents()  # M17: 166 -> (removed)
$ws.Cells.Item(19, 8).Value = 4999.8  # H19: 5000 -> 4999.8
$ws.Cells.Item(19, 9).Value = 4999  # I19: 0 -> 4999
$ws.Cells.Item(19, 11).Value = 14997  # K19: 0 -> 14997
$ws.Cells.Item(19, 13).Value = -14823  # M19: None -> -14823
$ws.Cells.Item(25, 8).Value = 1200  # H25: 0 -> 1200
$ws.Cells.Item(25, 10).Value = 1200  # J25: 0 -> 1200
$ws.Cells.Item(25, 12).Value = 3600  # L25: 0 -> 3600
$ws.Cells.Item(25, 14).Value = -3938  # N25: None -> -3938
$ws.Cells.Item(29, 8).Value = 55  # H29: 300 -> 55
$ws.Cells.Item(29, 9).Value = 100  # I29: 0 -> 100
$ws.Cells.Item(29, 10).Value = 10  # J29: 300 -> 10
$ws.Cells.Item(29, 11).Value = 300  # K29: 0 -> 300
$ws.Cells.Item(29, 12).Value = 30  # L29: 900 -> 30
$ws.Cells.Item(29, 13).Value = -23  # M29: None -> -23
$ws.Cells.Item(29, 14).Value = -584  # N29: -1454 -> -584
$ws.Cells.Item(30, 8).Value = 1200  # H30: 0 -> 1200
$ws.Cells.Item(30, 10).Value = 1200  # J30: 0 -> 1200
$ws.Cells.Item(30, 12).Value = 3600  # L30: 0 -> 3600
$ws.Cells.Item(30, 14).Value = -3804  # N30: None -> -3804
$ws.Cells.Item(34, 8).Value = 658.3333  # H34: 650 -> 658.3333
$ws.Cells.Item(34, 10).Value = 740  # J34: 750 -> 740
$ws.Cells.Item(34, 12).Value = 2220  # L34: 2250 -> 2220
$ws.Cells.Item(34, 14).Value = -2388  # N34: -2418 -> -2388
$ws.Cells.Item(44, 10).Value = 0  # J44: 1000 -> 0
$ws.Cells.Item(44, 12).Value = 0  # L44: 3000 -> 0
$ws.Cells.Item(44, 14).ClearContents()  # N44: -3796 -> (removed)
$ws.Cells.Item(46, 8).Value = 1200  # H46: 0 -> 1200
$ws.Cells.Item(46, 10).Value = 1200  # J46: 0 -> 1200
$ws.Cells.Item(46, 12).Value = 3600  # L46: 0 -> 3600
$ws.Cells.Item(46, 14).Value = -3782  # N46: None -> -3782
$ws.Cells.Item(58, 8).Value = 735.3333  # H58: 600 -> 735.3333
$ws.Cells.Item(58, 10).Value = 1006  # J58: 0 -> 1006
$ws.Cells.Item(58, 12).Value = 3018  # L58: 0 -> 3018
$ws.Cells.Item(58, 14).Value = -3274  # N58: None -> -3274
$ws.Cells.Item(63, 8).Value = 2514  # H63: 0 -> 2514
$ws.Cells.Item(63, 10).Value = 2514  # J63: 0 -> 2514
$ws.Cells.Item(63, 12).Value = 7542  # L63: 0 -> 7542
$ws.Cells.Item(63, 14).Value = -9040  # N63: None -> -9040
$ws.Cells.Item(64, 8).Value = 2856  # H64: 212 -> 2856
$ws.Cells.Item(64, 10).Value = 5500  # J64: 0 -> 5500
$ws.Cells.Item(64, 12).Value = 16500  # L64: 0 -> 16500
$ws.Cells.Item(64, 14).Value = -17040  # N64: None -> -17040
$ws.Cells.Item(66, 8).Value = 2514  # H66: 0 -> 2514
$ws.Cells.Item(66, 10).Value = 2514  # J66: 0 -> 2514
$ws.Cells.Item(66, 12).Value = 22626  # L66: 0 -> 22626
$ws.Cells.Item(66, 14).Value = -30114  # N66: None -> -30114
$ws.Cells.Item(67, 8).Value = 2856  # H67: 212 -> 2856
$ws.Cells.Item(67, 10).Value = 5500  # J67: 0 -> 5500
$ws.Cells.Item(67, 12).Value = 16500  # L67: 0 -> 16500
$ws.Cells.Item(67, 14).Value = -18372  # N67: None -> -18372
$ws.Cells.Item(75, 8).Value = 174  # H75: 198 -> 174
$ws.Cells.Item(75, 9).Value = 174  # I75: 198 -> 174
$ws.Cells.Item(75, 11).Value = 522  # K75: 594 -> 522
$ws.Cells.Item(75, 13).Value = 476  # M75: 404 -> 476
$ws.Cells.Item(78, 8).Value = 174  # H78: 198 -> 174
$ws.Cells.Item(78, 9).Value = 174  # I78: 198 -> 174
$ws.Cells.Item(78, 11).Value = 1566  # K78: 1782 -> 1566
$ws.Cells.Item(78, 13).Value = 3426  # M78: 3210 -> 3426

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 79.72727  # H2: 74.75 -> 79.72727
$ws.Cells.Item(2, 9).Value = 57.7  # I2: 54.272728 -> 57.7
$ws.Cells.Item(2, 11).Value = 57.7  # K2: 54.272728 -> 57.7
$ws.Cells.Item(2, 13).Value = 55.3  # M2: 58.727272 -> 55.3
$ws.Cells.Item(55, 8).Value = 8000  # H55: 10000 -> 8000
$ws.Cells.Item(55, 10).Value = 8000  # J55: 10000 -> 8000
$ws.Cells.Item(55, 12).Value = 8000  # L55: 10000 -> 8000
$ws.Cells.Item(55, 14).Value = -8654  # N55: -10654 -> -8654

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(63, 8).Value = 40000  # H63: 49995 -> 40000
$ws.Cells.Item(63, 10).Value = 40000  # J63: 49995 -> 40000
$ws.Cells.Item(63, 12).Value = 40000  # L63: 49995 -> 40000
$ws.Cells.Item(63, 14).Value = -41248  # N63: -51243 -> -41248
$ws.Cells.Item(66, 8).Value = 40000  # H66: 49995 -> 40000
$ws.Cells.Item(66, 10).Value = 40000  # J66: 49995 -> 40000
$ws.Cells.Item(66, 12).Value = 120000  # L66: 149985 -> 120000
$ws.Cells.Item(66, 14).Value = -126240  # N66: -156225 -> -126240
$ws.Cells.Item(68, 8).Value = 26774  # H68: 0 -> 26774
$ws.Cells.Item(68, 10).Value = 26774  # J68: 0 -> 26774
$ws.Cells.Item(68, 12).Value = 26774  # L68: 0 -> 26774
$ws.Cells.Item(68, 14).Value = -28396  # N68: None -> -28396
$ws.Cells.Item(69, 8).Value = 30000  # H69: 27499.75 -> 30000
$ws.Cells.Item(69, 10).Value = 30000  # J69: 27499.75 -> 30000
$ws.Cells.Item(69, 12).Value = 30000  # L69: 27499.75 -> 30000
$ws.Cells.Item(69, 14).Value = -31498  # N69: -28997.75 -> -31498
$ws.Cells.Item(71, 8).Value = 26774  # H71: 0 -> 26774
$ws.Cells.Item(71, 10).Value = 26774  # J71: 0 -> 26774
$ws.Cells.Item(71, 12).Value = 80322  # L71: 0 -> 80322
$ws.Cells.Item(71, 14).Value = -88434  # N71: None -> -88434
$ws.Cells.Item(72, 8).Value = 30000  # H72: 27499.75 -> 30000
$ws.Cells.Item(72, 10).Value = 30000  # J72: 27499.75 -> 30000
$ws.Cells.Item(72, 12).Value = 90000  # L72: 82499.25 -> 90000
$ws.Cells.Item(72, 14).Value = -97488  # N72: -89987.25 -> -97488
$ws.Cells.Item(126, 8).Value = 2299.75  # H126: 2666.3333 -> 2299.75
$ws.Cells.Item(126, 9).Value = 1399.6666  # I126: 1499.5 -> 1399.6666
$ws.Cells.Item(126, 11).Value = 4198.9998  # K126: 4498.5 -> 4198.9998
$ws.Cells.Item(126, 13).Value = -1728.9998  # M126: -2028.5 -> -1728.9998
